# Updated cryptos list values (Price and Volume(1h) columns) per data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.832.00"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "'2.241.52"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'112.57"
$ws.Range("E5").Value = "  -8.47%  "
$ws.Range("D6").Value = "'295.80"
$ws.Range("E6").Value = "  +10.46%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").Value = "'45.65"
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("D11").Value = "'0.0925"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Value = "'55.74"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "'9.00"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'0.104"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").Value = "'0.912"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "'15.26"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "'2.580.89"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "'2.262.26"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "'42.716.00"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").Value = "'7.51"
$ws.Range("E20").Value = "  +5.35%  "
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("D22").Value = "'73.28"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "'3.55"
$ws.Range("E23").Value = "  +22.17%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("D25").Value = "'231.41"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").Value = "'9.43"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "'12.00"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'39.94"
$ws.Range("E29").Value = "  -6.59%  "
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("D32").Value = "'174.12"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'21.25"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "'0.0897"
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("D35").Value = "'5.75"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").Value = "'5.05"
$ws.Range("E36").Value = "  +6.25%  "
$ws.Range("D37").Value = "'4.30"
$ws.Range("E37").Value = "  +6.65%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'0.0371"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").Value = "'2.55"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'0.240"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'71.78"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").Value = "'13.26"
$ws.Range("E44").Value = "  -7.80%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "'1.33"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("E47").Value = "  -6.73%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").Value = "'106.30"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("D50").Value = "'8.67"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.0989"
$ws.Range("E51").Value = "  -1.80%  "
